$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.762.02'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '2.897.69'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'197.42"
$ws.Range("E5").Value = '  +4.12%  '
$ws.Range("D6").Value = "'595.01"
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = "'0.549"
$ws.Range("E8").Value = '  -3.71%  '
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").Value = '2.895.28'
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("D11").Value = "'0.418"
$ws.Range("E11").Value = '  +10.97%  '
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("E13").Value = '  -2.70%  '
$ws.Range("D14").Value = '3.428.83'
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").Value = '75.707.10'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").Value = "'0.0000188"
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("D17").Value = "'27.22"
$ws.Range("E17").Value = '  -2.38%  '
$ws.Range("D18").Value = '2.897.39'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("E19").Value = '  -4.58%  '
$ws.Range("D20").Value = "'12.54"
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("D21").Value = "'375.65"
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = "'4.14"
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").Value = "'70.94"
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("D26").Value = '3.038.22'
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").Value = "'4.19"
$ws.Range("E27").Value = '  -2.79%  '
$ws.Range("D28").Value = "'9.52"
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("E29").Value = '  +1.55%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("E31").Value = '  -2.88%  '
$ws.Range("D32").Value = "'501.61"
$ws.Range("E32").Value = '  -6.53%  '
$ws.Range("D33").Value = "'7.67"
$ws.Range("E33").Value = '  -4.23%  '
$ws.Range("D34").Value = "'1.79"
$ws.Range("E34").Value = '  -2.67%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = "'163.40"
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").Value = "'19.94"
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("E39").Value = '  -7.39%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = "'179.03"
$ws.Range("E41").Value = '  -3.17%  '
$ws.Range("D42").Value = "'0.340"
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("D43").Value = "'4.98"
$ws.Range("E43").Value = '  -4.29%  '
$ws.Range("D44").Value = "'1.65"
$ws.Range("E44").Value = '  -3.55%  '
$ws.Range("D45").Value = "'0.0901"
$ws.Range("E45").Value = '  +4.53%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = "'40.03"
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("B47").Value = 'ImmutableX'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = '  -5.70%  '
$ws.Range("E48").Value = '  -2.99%  '
$ws.Range("D49").Value = "'0.573"
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("E50").Value = '  -2.72%  '
$ws.Range("D51").Value = "'0.649"
$ws.Range("E51").Value = '  +4.81%  '
